$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Bank Angle"
$ws.Range("B8").Value = "Theta_bank"
$ws.Range("C8").Value = "deg"

$ws.Range("A8:C8").Select()
